# Auto-generated Excel COM-interop script
# Applies the 'Update countries & provincias Spain' data refresh:
#  - refreshes the COVID country statistics table (sorted desc by Casos totales)
#  - updates the 'Datos actualizados' timestamp in A1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 9 de Octubre de 2020 a las 15:42'
$ws.Cells.Item(4, 2).Value = 7838537
$ws.Cells.Item(4, 3).Value = 4774
$ws.Cells.Item(4, 4).Value = 5026289
$ws.Cells.Item(4, 5).Value = 2594431
$ws.Cells.Item(4, 7).Value = 79
$ws.Cells.Item(4, 8).Value = 217817
$ws.Cells.Item(5, 2).Value = 6908603
$ws.Cells.Item(5, 3).Value = 4791
$ws.Cells.Item(5, 5).Value = 895982
$ws.Cells.Item(5, 7).Value = 31
$ws.Cells.Item(5, 8).Value = 106552
$ws.Cells.Item(16, 2).Value = 492378
$ws.Cells.Item(16, 3).Value = 4142
$ws.Cells.Item(16, 4).Value = 401379
$ws.Cells.Item(16, 5).Value = 62901
$ws.Cells.Item(16, 7).Value = 210
$ws.Cells.Item(16, 8).Value = 28098
$ws.Cells.Item(18, 2).Value = 397780
$ws.Cells.Item(18, 3).Value = 3214
$ws.Cells.Item(18, 4).Value = 328097
$ws.Cells.Item(18, 5).Value = 59948
$ws.Cells.Item(18, 7).Value = 52
$ws.Cells.Item(18, 8).Value = 9735
$ws.Cells.Item(20, 1).Value = 'Arabia Saudita'
$ws.Cells.Item(20, 2).Value = 338539
$ws.Cells.Item(20, 3).Value = 407
$ws.Cells.Item(20, 4).Value = 324282
$ws.Cells.Item(20, 5).Value = 9261
$ws.Cells.Item(20, 7).Value = 24
$ws.Cells.Item(20, 8).Value = 4996
$ws.Cells.Item(21, 1).Value = 'Italia'
$ws.Cells.Item(21, 2).Value = 338398
$ws.Cells.Item(21, 4).Value = 236363
$ws.Cells.Item(21, 5).Value = 65952
$ws.Cells.Item(21, 8).Value = 36083
$ws.Cells.Item(26, 2).Value = 316859
$ws.Cells.Item(26, 3).Value = 1345
$ws.Cells.Item(26, 5).Value = 37689
$ws.Cells.Item(30, 2).Value = 161781
$ws.Cells.Item(30, 3).Value = 5971
$ws.Cells.Item(30, 7).Value = 13
$ws.Cells.Item(30, 8).Value = 6544
$ws.Cells.Item(40, 2).Value = 110076
$ws.Cells.Item(40, 3).Value = 635
$ws.Cells.Item(40, 4).Value = 102024
$ws.Cells.Item(40, 5).Value = 7403
$ws.Cells.Item(40, 7).Value = 7
$ws.Cells.Item(40, 8).Value = 649
$ws.Cells.Item(52, 2).Value = 83928
$ws.Cells.Item(52, 3).Value = 1394
$ws.Cells.Item(52, 4).Value = 52164
$ws.Cells.Item(52, 5).Value = 29702
$ws.Cells.Item(52, 7).Value = 12
$ws.Cells.Item(52, 8).Value = 2062
$ws.Cells.Item(70, 2).Value = 43664
$ws.Cells.Item(70, 3).Value = 408
$ws.Cells.Item(70, 4).Value = 36922
$ws.Cells.Item(70, 5).Value = 6375
$ws.Cells.Item(70, 7).Value = 8
$ws.Cells.Item(70, 8).Value = 367
$ws.Cells.Item(72, 2).Value = 41368
$ws.Cells.Item(72, 3).Value = 1076
$ws.Cells.Item(72, 4).Value = 23453
$ws.Cells.Item(72, 5).Value = 17294
$ws.Cells.Item(72, 7).Value = 5
$ws.Cells.Item(72, 8).Value = 621
$ws.Cells.Item(77, 2).Value = 34517
$ws.Cells.Item(77, 3).Value = 173
$ws.Cells.Item(77, 5).Value = 2220
$ws.Cells.Item(77, 8).Value = 761
$ws.Cells.Item(79, 2).Value = 29951
$ws.Cells.Item(79, 3).Value = 109
$ws.Cells.Item(79, 4).Value = 24995
$ws.Cells.Item(79, 5).Value = 4075
$ws.Cells.Item(80, 2).Value = 29917
$ws.Cells.Item(80, 3).Value = 389
$ws.Cells.Item(80, 4).Value = 23241
$ws.Cells.Item(80, 5).Value = 5750
$ws.Cells.Item(80, 7).Value = 13
$ws.Cells.Item(80, 8).Value = 926
$ws.Cells.Item(84, 1).Value = 'Birmania'
$ws.Cells.Item(84, 2).Value = 23906
$ws.Cells.Item(84, 3).Value = 1461
$ws.Cells.Item(84, 4).Value = 6738
$ws.Cells.Item(84, 5).Value = 16602
$ws.Cells.Item(84, 7).Value = 31
$ws.Cells.Item(84, 8).Value = 566
$ws.Cells.Item(85, 1).Value = 'Bulgaria'
$ws.Cells.Item(85, 2).Value = 23259
$ws.Cells.Item(85, 4).Value = 15563
$ws.Cells.Item(85, 5).Value = 6816
$ws.Cells.Item(85, 8).Value = 880
$ws.Cells.Item(89, 1).Value = 'Republica de Macedonia'
$ws.Cells.Item(89, 2).Value = 20163
$ws.Cells.Item(89, 3).Value = 386
$ws.Cells.Item(89, 4).Value = 16009
$ws.Cells.Item(89, 5).Value = 3373
$ws.Cells.Item(89, 7).Value = 6
$ws.Cells.Item(89, 8).Value = 781
$ws.Cells.Item(90, 1).Value = 'Costa de Marfil'
$ws.Cells.Item(90, 2).Value = 19982
$ws.Cells.Item(90, 4).Value = 19626
$ws.Cells.Item(90, 5).Value = 236
$ws.Cells.Item(90, 8).Value = 120
$ws.Cells.Item(93, 2).Value = 16676
$ws.Cells.Item(93, 3).Value = 22
$ws.Cells.Item(93, 4).Value = 15975
$ws.Cells.Item(93, 5).Value = 464
$ws.Cells.Item(93, 7).Value = 2
$ws.Cells.Item(93, 8).Value = 237
$ws.Cells.Item(107, 1).Value = 'Tayikistan'
$ws.Cells.Item(107, 2).Value = 10137
$ws.Cells.Item(107, 3).Value = 40
$ws.Cells.Item(107, 4).Value = 8959
$ws.Cells.Item(107, 5).Value = 1099
$ws.Cells.Item(107, 7).Value = 1
$ws.Cells.Item(107, 8).Value = 79
$ws.Cells.Item(108, 1).Value = 'Guayana Francesa'
$ws.Cells.Item(108, 2).Value = 10128
$ws.Cells.Item(108, 4).Value = 9799
$ws.Cells.Item(108, 5).Value = 260
$ws.Cells.Item(108, 8).Value = 69

Write-Host "Applied 123 cell updates"
